# Fix steam reforming elc/lng ratio (column C = "Correct CF") for the
# 8 timeslices on the New_COMFR sheet. Dependent formulas (F, I columns,
# the row-12 sums, and the F16:F21 / M17:M20 ratio columns) recalculate
# automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New_COMFR")
$ws.Activate()

$ws.Range("C4").Value = 0.06
$ws.Range("C5").Value = 0.2
$ws.Range("C6").Value = 0.04
$ws.Range("C7").Value = 0.2
$ws.Range("C8").Value = 0.06
$ws.Range("C9").Value = 0.2
$ws.Range("C10").Value = 0.04
$ws.Range("C11").Value = 0.2

# Leave the selection on the edited range, matching the saved UI state.
$ws.Range("C4:C11").Select() | Out-Null
